$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q3" sheet so the OLD fund-level data
#    survives unchanged in a new tab placed right after it. This new
#    tab keeps the old formatting/margins/styles untouched and will be
#    renamed to "2022-Q3" at the end (taking over that name).
# ------------------------------------------------------------------
$wsOldQ3 = $wb.Worksheets.Item("2022-Q3")
$wsOldQ3.Copy($null, $wsOldQ3)

$wsDuplicate = $wb.Worksheets.Item("2022-Q3 (2)")
$wsDuplicate.Name = "2022-Q3-STAGING"

# The original tab (still named "2022-Q3", still holding the OLD data)
# becomes the new "2022-Q4" sheet - we keep its identity (sheetId/rId)
# but swap in the new quarter's data below.
$wsQ4 = $wb.Worksheets.Item("2022-Q3")
$wsQ4.Name = "2022-Q4"

# The staged duplicate becomes the final "2022-Q3" tab - its content is
# already correct (untouched old data), nothing else to do for it.
$wsDuplicate.Name = "2022-Q3"
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# ------------------------------------------------------------------
# 2. Clear out the old data on the (renamed) "2022-Q4" sheet and fill
#    it in with the new quarter's fund holdings.
# ------------------------------------------------------------------
$wsQ4.Cells.Clear()

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("C2").Value = "汇添富全球移动互联灵活配置混合（QDII）D"
$wsQ4.Range("H2").Value = 8

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("C3").Value = "汇添富全球移动互联灵活配置混合（QDII）A"
$wsQ4.Range("H3").Value = 8

$wsQ4.Range("A4").Value = 2
$wsQ4.Range("C4").Value = "汇添富全球移动互联灵活配置混合（QDII）C"
$wsQ4.Range("H4").Value = 8

# The fund codes / scale / position figures are stored as literal TEXT
# in this workbook (not numbers), even though they look numeric (e.g.
# "015203" with a leading zero). Force the cells to Text format before
# typing them in so Excel does not silently convert them to numbers.
$textCells = $wsQ4.Range("B2:B4,D2:G4")
$textCells.NumberFormat = "@"

$wsQ4.Range("B2").Value = "015203"
$wsQ4.Range("D2").Value = "11.52"
$wsQ4.Range("E2").Value = "92.14"
$wsQ4.Range("F2").Value = "3.22"
$wsQ4.Range("G2").Value = "0.3709"

$wsQ4.Range("B3").Value = "001668"
$wsQ4.Range("D3").Value = "11.48"
$wsQ4.Range("E3").Value = "92.14"
$wsQ4.Range("F3").Value = "3.22"
$wsQ4.Range("G3").Value = "0.3697"

$wsQ4.Range("B4").Value = "015202"
$wsQ4.Range("D4").Value = "0.01"
$wsQ4.Range("E4").Value = "92.14"
$wsQ4.Range("F4").Value = "3.22"
$wsQ4.Range("G4").Value = "0.0003"

# Drop the temporary "@" number format back off those cells (an unused
# cell never touched carries the workbook's plain default formatting)
# so the cells end up as plain text values with no special styling.
$wsQ4.Range("ZZ1").Copy() | Out-Null
$textCells.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the header / row-number formatting used throughout this
# workbook (same style as the "总计" sheet's header row) by copying it
# over, rather than building a brand-new style.
$ws1 = $wb.Worksheets.Item("总计")
$ws1.Range("B1").Copy() | Out-Null
$wsQ4.Range("B1:H1").PasteSpecial(-4122)

$ws1.Range("A2").Copy() | Out-Null
$wsQ4.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the "总计" sheet's page margins (0.75in/1in/0.5in) on the new
# "2022-Q4" tab.
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 3. Update the "总计" summary sheet: push the 2022-Q3 summary row
#    down to row 3 and add the new 2022-Q4 summary row in row 2.
# ------------------------------------------------------------------
$ws1.Range("A2").Copy() | Out-Null
$ws1.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q3"
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 0.32

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 0.74

# ------------------------------------------------------------------
# 4. Make sure the "总计" sheet stays the active one, matching the
#    original workbook's selection state.
# ------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
